# fixing matricula of Matc65
# Update column A (matricula) values for rows 14-39 according to the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    14 = "217216526"
    15 = "216117974"
    16 = "221117463"
    17 = "217125254"
    18 = "219218129"
    19 = "218215397"
    20 = "220117282"
    21 = "219217429"
    22 = "216216087"
    23 = "220121412"
    24 = "210201260"
    25 = "201520233"
    26 = "217117994"
    27 = "219118481"
    28 = "221119218"
    29 = "219215012"
    30 = "219121541"
    31 = "214007731"
    32 = "219215013"
    33 = "220117290"
    34 = "219118473"
    35 = "220117273"
    36 = "220120071"
    37 = "221216783"
    38 = "214120645"
    39 = "220217140"
}

foreach ($row in $newValues.Keys) {
    # Prefix with an apostrophe so Excel stores the numeric-looking
    # matricula id as text (matching the original inlineStr/text cells)
    # instead of auto-converting it to a number.
    $ws.Cells.Item($row, 1).Value = "'" + $newValues[$row]
}
